$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 657, pushing the existing
# rows 657-692 down to become rows 659-694.
$ws.Rows.Item(657).Insert()
$ws.Rows.Item(657).Insert()

# New row 657: "1a plateado" entry for 2023-01-13
$ws.Cells.Item(657, 1).Value = 4
$ws.Cells.Item(657, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(657, 3).Value = "Los Lagos"
$ws.Cells.Item(657, 4).Value = 44939
$ws.Cells.Item(657, 5).Value = 10
$ws.Cells.Item(657, 6).Value = "Fruta"
$ws.Cells.Item(657, 7).Value = 100102
$ws.Cells.Item(657, 8).Value = "Cítricos"
$ws.Cells.Item(657, 9).Value = 100102003
$ws.Cells.Item(657, 10).Value = "Limón"
$ws.Cells.Item(657, 11).Value = "Sin especificar"
$ws.Cells.Item(657, 12).Value = "1a plateado"
$ws.Cells.Item(657, 13).Value = 1200
$ws.Cells.Item(657, 14).Value = 21000
$ws.Cells.Item(657, 15).Value = 22000
$ws.Cells.Item(657, 16).Value = 21500
$ws.Cells.Item(657, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(657, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(657, 19).Value = 1194
$ws.Cells.Item(657, 20).Value = 18

# New row 658: "2a plateado" entry for 2023-01-13
$ws.Cells.Item(658, 1).Value = 4
$ws.Cells.Item(658, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(658, 3).Value = "Los Lagos"
$ws.Cells.Item(658, 4).Value = 44939
$ws.Cells.Item(658, 5).Value = 10
$ws.Cells.Item(658, 6).Value = "Fruta"
$ws.Cells.Item(658, 7).Value = 100102
$ws.Cells.Item(658, 8).Value = "Cítricos"
$ws.Cells.Item(658, 9).Value = 100102003
$ws.Cells.Item(658, 10).Value = "Limón"
$ws.Cells.Item(658, 11).Value = "Sin especificar"
$ws.Cells.Item(658, 12).Value = "2a plateado"
$ws.Cells.Item(658, 13).Value = 600
$ws.Cells.Item(658, 14).Value = 19000
$ws.Cells.Item(658, 15).Value = 19000
$ws.Cells.Item(658, 16).Value = 19000
$ws.Cells.Item(658, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(658, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(658, 19).Value = 1056
$ws.Cells.Item(658, 20).Value = 18
